$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 260.8889
$ws.Cells.Item(2, 9).Value = 98.2
$ws.Cells.Item(2, 10).Value = 464.25
$ws.Cells.Item(2, 11).Value = 98.2
$ws.Cells.Item(2, 12).Value = 464.25
$ws.Cells.Item(2, 13).Value = 14.8
$ws.Cells.Item(2, 14).Value = -690.25
$ws.Cells.Item(17, 8).Value = 1472604.6
$ws.Cells.Item(17, 10).Value = 1517208.4
$ws.Cells.Item(17, 12).Value = 4551625.199999999
$ws.Cells.Item(17, 14).Value = -4551961.199999999
$ws.Cells.Item(43, 8).Value = 1292.6666
$ws.Cells.Item(43, 9).Value = 950
$ws.Cells.Item(43, 10).Value = 1978
$ws.Cells.Item(43, 11).Value = 950
$ws.Cells.Item(43, 12).Value = 1978
$ws.Cells.Item(43, 13).Value = -881
$ws.Cells.Item(43, 14).Value = -2116
$ws.Cells.Item(135, 8).Value = 27785870
$ws.Cells.Item(135, 9).Value = 761.0714
$ws.Cells.Item(135, 10).Value = 125033750
$ws.Cells.Item(135, 11).Value = 6849.6426
$ws.Cells.Item(135, 12).Value = 1125303750
$ws.Cells.Item(135, 13).Value = -4314.6426
$ws.Cells.Item(135, 14).Value = -1125308820
$ws.Cells.Item(137, 8).Value = 2827.7727
$ws.Cells.Item(137, 9).Value = 2661.7222
$ws.Cells.Item(137, 10).Value = 3575
$ws.Cells.Item(137, 11).Value = 7985.1666
$ws.Cells.Item(137, 12).Value = 10725
$ws.Cells.Item(137, 13).Value = -5435.1666
$ws.Cells.Item(137, 14).Value = -15825

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 2833.3333
$ws.Cells.Item(22, 9).Value = 1750
$ws.Cells.Item(22, 10).Value = 5000
$ws.Cells.Item(22, 11).Value = 1750
$ws.Cells.Item(22, 12).Value = 5000
$ws.Cells.Item(22, 13).Value = -1451
$ws.Cells.Item(22, 14).Value = -5598
$ws.Cells.Item(32, 8).Value = 4675.269
$ws.Cells.Item(32, 9).Value = 4642
$ws.Cells.Item(32, 10).Value = 5507
$ws.Cells.Item(32, 11).Value = 4642
$ws.Cells.Item(32, 12).Value = 5507
$ws.Cells.Item(32, 13).Value = -4355
$ws.Cells.Item(32, 14).Value = -6081
$ws.Cells.Item(45, 8).Value = 1977.1875
$ws.Cells.Item(45, 10).Value = 1695.7858
$ws.Cells.Item(45, 12).Value = 1695.7858
$ws.Cells.Item(45, 14).Value = -2449.7858
$ws.Cells.Item(61, 8).Value = 2052.3823
$ws.Cells.Item(61, 9).Value = 1621.862
$ws.Cells.Item(61, 11).Value = 1621.862
$ws.Cells.Item(61, 13).Value = -1409.862
$ws.Cells.Item(74, 8).Value = 125001680
$ws.Cells.Item(74, 9).Value = 250000720
$ws.Cells.Item(74, 11).Value = 250000720
$ws.Cells.Item(74, 13).Value = -249999846
$ws.Cells.Item(77, 8).Value = 125001680
$ws.Cells.Item(77, 9).Value = 250000720
$ws.Cells.Item(77, 11).Value = 1250003600
$ws.Cells.Item(77, 13).Value = -1249999232
$ws.Cells.Item(102, 8).Value = 1747.375
$ws.Cells.Item(102, 9).Value = 1245
$ws.Cells.Item(102, 11).Value = 1245
$ws.Cells.Item(102, 13).Value = 377
$ws.Cells.Item(110, 8).Value = 1333.3334
$ws.Cells.Item(110, 9).Value = 1000
$ws.Cells.Item(110, 10).Value = 1500
$ws.Cells.Item(110, 11).Value = 1000
$ws.Cells.Item(110, 12).Value = 1500
$ws.Cells.Item(110, 13).Value = 1045
$ws.Cells.Item(110, 14).Value = -5590
$ws.Cells.Item(122, 8).Value = 2256.4614
$ws.Cells.Item(122, 9).Value = 1592
$ws.Cells.Item(122, 10).Value = 3319.6
$ws.Cells.Item(122, 11).Value = 4776
$ws.Cells.Item(122, 12).Value = 9958.799999999999
$ws.Cells.Item(122, 13).Value = -2326
$ws.Cells.Item(122, 14).Value = -14858.8
$ws.Cells.Item(132, 8).Value = 14688.257
$ws.Cells.Item(132, 9).Value = 1741
$ws.Cells.Item(132, 11).Value = 5223
$ws.Cells.Item(132, 13).Value = -2693
$ws.Cells.Item(136, 8).Value = 2052.3823
$ws.Cells.Item(136, 9).Value = 1621.862
$ws.Cells.Item(136, 11).Value = 4865.586
$ws.Cells.Item(136, 13).Value = -2315.586

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2163.5454
$ws.Cells.Item(99, 9).Value = 1759.8
$ws.Cells.Item(99, 11).Value = 1759.8
$ws.Cells.Item(99, 13).Value = -261.8
$ws.Cells.Item(134, 8).Value = 4047.7144
$ws.Cells.Item(134, 9).Value = 4243.6924
$ws.Cells.Item(134, 10).Value = 1500
$ws.Cells.Item(134, 11).Value = 12731.0772
$ws.Cells.Item(134, 12).Value = 4500
$ws.Cells.Item(134, 13).Value = -10196.0772
$ws.Cells.Item(134, 14).Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 11829.409
$ws.Cells.Item(31, 9).Value = 16696.154
$ws.Cells.Item(31, 10).Value = 4799.6665
$ws.Cells.Item(31, 11).Value = 16696.154
$ws.Cells.Item(31, 12).Value = 4799.6665
$ws.Cells.Item(31, 13).Value = -16401.154
$ws.Cells.Item(31, 14).Value = -5389.6665
$ws.Cells.Item(34, 8).Value = 11829.409
$ws.Cells.Item(34, 9).Value = 16696.154
$ws.Cells.Item(34, 10).Value = 4799.6665
$ws.Cells.Item(34, 11).Value = 16696.154
$ws.Cells.Item(34, 12).Value = 4799.6665
$ws.Cells.Item(34, 13).Value = -16494.154
$ws.Cells.Item(34, 14).Value = -5203.6665
$ws.Cells.Item(58, 8).Value = 15800
$ws.Cells.Item(58, 9).Value = 1135.4348
$ws.Cells.Item(58, 10).Value = 46462.273
$ws.Cells.Item(58, 11).Value = 1135.4348
$ws.Cells.Item(58, 12).Value = 46462.273
$ws.Cells.Item(58, 13).Value = -932.4348
$ws.Cells.Item(58, 14).Value = -46868.273
$ws.Cells.Item(86, 8).Value = 10428147
$ws.Cells.Item(86, 9).Value = 2646.6667
$ws.Cells.Item(86, 10).Value = 23832362
$ws.Cells.Item(86, 11).Value = 2646.6667
$ws.Cells.Item(86, 12).Value = 23832362
$ws.Cells.Item(86, 13).Value = -1523.6667
$ws.Cells.Item(86, 14).Value = -23834608
$ws.Cells.Item(89, 8).Value = 10428147
$ws.Cells.Item(89, 9).Value = 2646.6667
$ws.Cells.Item(89, 10).Value = 23832362
$ws.Cells.Item(89, 11).Value = 13233.3335
$ws.Cells.Item(89, 12).Value = 119161810
$ws.Cells.Item(89, 13).Value = -7617.333500000001
$ws.Cells.Item(89, 14).Value = -119173042
$ws.Cells.Item(122, 8).Value = 1093.25
$ws.Cells.Item(122, 9).Value = 1186.625
$ws.Cells.Item(122, 10).Value = 999.875
$ws.Cells.Item(122, 11).Value = 3559.875
$ws.Cells.Item(122, 12).Value = 2999.625
$ws.Cells.Item(122, 13).Value = -1109.875
$ws.Cells.Item(122, 14).Value = -7899.625
$ws.Cells.Item(132, 8).Value = 12218.58
$ws.Cells.Item(132, 9).Value = 15516.944
$ws.Cells.Item(132, 10).Value = 3737.0715
$ws.Cells.Item(132, 11).Value = 46550.83199999999
$ws.Cells.Item(132, 12).Value = 11211.2145
$ws.Cells.Item(132, 13).Value = -44020.83199999999
$ws.Cells.Item(132, 14).Value = -16271.2145
$ws.Cells.Item(134, 8).Value = 1271.8214
$ws.Cells.Item(134, 9).Value = 897.2941
$ws.Cells.Item(134, 11).Value = 2691.8823
$ws.Cells.Item(134, 13).Value = -156.8822999999998
$ws.Cells.Item(136, 8).Value = 15800
$ws.Cells.Item(136, 9).Value = 1135.4348
$ws.Cells.Item(136, 10).Value = 46462.273
$ws.Cells.Item(136, 11).Value = 3406.3044
$ws.Cells.Item(136, 12).Value = 139386.819
$ws.Cells.Item(136, 13).Value = -856.3044
$ws.Cells.Item(136, 14).Value = -144486.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 697.4
$ws.Cells.Item(5, 9).Value = 700.8889
$ws.Cells.Item(5, 10).Value = 666
$ws.Cells.Item(5, 11).Value = 2102.6667
$ws.Cells.Item(5, 12).Value = 1998
$ws.Cells.Item(5, 13).Value = -1990.6667
$ws.Cells.Item(5, 14).Value = -2222
$ws.Cells.Item(47, 8).Value = 88.25
$ws.Cells.Item(47, 9).Value = 88.25
$ws.Cells.Item(47, 11).Value = 264.75
$ws.Cells.Item(47, 13).Value = 166.25
$ws.Cells.Item(122, 8).Value = 1421.5714
$ws.Cells.Item(122, 10).Value = 1530.6842
$ws.Cells.Item(122, 12).Value = 13776.1578
$ws.Cells.Item(122, 14).Value = -18676.1578
$ws.Cells.Item(129, 8).Value = 295508.25
$ws.Cells.Item(129, 9).Value = 718.2
$ws.Cells.Item(129, 10).Value = 418337.4
$ws.Cells.Item(129, 11).Value = 2154.6
$ws.Cells.Item(129, 12).Value = 1255012.2
$ws.Cells.Item(129, 13).Value = 2845.4
$ws.Cells.Item(129, 14).Value = -1265012.2
$ws.Cells.Item(131, 8).Value = 786.47
$ws.Cells.Item(131, 10).Value = 786.47
$ws.Cells.Item(131, 12).Value = 2359.41
$ws.Cells.Item(131, 14).Value = -12439.41
$ws.Cells.Item(135, 8).Value = 697.4
$ws.Cells.Item(135, 9).Value = 700.8889
$ws.Cells.Item(135, 10).Value = 666
$ws.Cells.Item(135, 11).Value = 6308.0001
$ws.Cells.Item(135, 12).Value = 5994
$ws.Cells.Item(135, 13).Value = -3773.0001
$ws.Cells.Item(135, 14).Value = -11064

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 49383804
$ws.Cells.Item(122, 9).Value = 19609086
$ws.Cells.Item(122, 10).Value = 100000830
$ws.Cells.Item(122, 11).Value = 58827258
$ws.Cells.Item(122, 12).Value = 300002490
$ws.Cells.Item(122, 13).Value = -58824808
$ws.Cells.Item(122, 14).Value = -300007390
$ws.Cells.Item(132, 8).Value = 28265.428
$ws.Cells.Item(132, 9).Value = 4651.7334
$ws.Cells.Item(132, 10).Value = 87299.664
$ws.Cells.Item(132, 11).Value = 13955.2002
$ws.Cells.Item(132, 12).Value = 261898.992
$ws.Cells.Item(132, 13).Value = -11425.2002
$ws.Cells.Item(132, 14).Value = -266958.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2885.2258
$ws.Cells.Item(7, 9).Value = 3982.9092
$ws.Cells.Item(7, 11).Value = 3982.9092
$ws.Cells.Item(7, 13).Value = -3870.9092
$ws.Cells.Item(100, 8).Value = 2764.2856
$ws.Cells.Item(100, 9).Value = 2340.6
$ws.Cells.Item(100, 10).Value = 2999.6667
$ws.Cells.Item(100, 11).Value = 2340.6
$ws.Cells.Item(100, 12).Value = 2999.6667
$ws.Cells.Item(100, 13).Value = -1799.6
$ws.Cells.Item(100, 14).Value = -4081.6667
$ws.Cells.Item(126, 8).Value = 2885.2258
$ws.Cells.Item(126, 9).Value = 3982.9092
$ws.Cells.Item(126, 11).Value = 11948.7276
$ws.Cells.Item(126, 13).Value = -9478.7276
$ws.Cells.Item(132, 8).Value = 2616.4517
$ws.Cells.Item(132, 9).Value = 1271
$ws.Cells.Item(132, 10).Value = 4250.2144
$ws.Cells.Item(132, 11).Value = 3813
$ws.Cells.Item(132, 12).Value = 12750.6432
$ws.Cells.Item(132, 13).Value = -1283
$ws.Cells.Item(132, 14).Value = -17810.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1540
$ws.Cells.Item(96, 9).Value = 1733.3334
$ws.Cells.Item(96, 10).Value = 1250
$ws.Cells.Item(96, 11).Value = 1733.3334
$ws.Cells.Item(96, 12).Value = 1250
$ws.Cells.Item(96, 13).Value = -360.3334
$ws.Cells.Item(96, 14).Value = -3996
$ws.Cells.Item(107, 8).Value = 4132909.2
$ws.Cells.Item(107, 9).Value = 916.6667
$ws.Cells.Item(107, 10).Value = 9091301
$ws.Cells.Item(107, 11).Value = 2750.0001
$ws.Cells.Item(107, 12).Value = 27273903
$ws.Cells.Item(107, 13).Value = -830.0001000000002
$ws.Cells.Item(107, 14).Value = -27277743
$ws.Cells.Item(132, 8).Value = 1349.909
$ws.Cells.Item(132, 9).Value = 927.7778
$ws.Cells.Item(132, 10).Value = 3249.5
$ws.Cells.Item(132, 11).Value = 2783.3334
$ws.Cells.Item(132, 12).Value = 9748.5
$ws.Cells.Item(132, 13).Value = -253.3334
$ws.Cells.Item(132, 14).Value = -14808.5
$ws.Cells.Item(136, 8).Value = 40002280
$ws.Cells.Item(136, 9).Value = 76925380
$ws.Cells.Item(136, 10).Value = 2258.6667
$ws.Cells.Item(136, 11).Value = 230776140
$ws.Cells.Item(136, 12).Value = 6776.000100000001
$ws.Cells.Item(136, 13).Value = -230773590

Write-Host "Applied Typhon_Profits updates"